$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.722.70"
$ws.Range("E2").Value = "  -3.81%  "

$ws.Range("D3").Value = "3.147.12"
$ws.Range("E3").Value = "  -4.10%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "215.61"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").Value = "634.02"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").Value = "0.399"
$ws.Range("E7").Value = "  -3.82%  "

$ws.Range("D8").Value = "0.731"
$ws.Range("E8").Value = "  +2.74%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").Value = "3.146.58"
$ws.Range("E10").Value = "  -4.04%  "

$ws.Range("D11").Value = "0.556"
$ws.Range("E11").Value = "  -5.67%  "

$ws.Range("E12").Value = "  -0.47%  "

$ws.Range("D13").Value = "0.0000252"
$ws.Range("E13").Value = "  -5.07%  "

$ws.Range("D14").Value = "5.30"
$ws.Range("E14").Value = "  -1.37%  "

$ws.Range("D15").Value = "88.565.88"
$ws.Range("E15").Value = "  -3.51%  "

$ws.Range("D16").Value = "3.708.24"
$ws.Range("E16").Value = "  -4.61%  "

$ws.Range("D17").Value = "32.55"
$ws.Range("E17").Value = "  -4.86%  "

$ws.Range("D18").Value = "3.140.21"
$ws.Range("E18").Value = "  -4.95%  "

$ws.Range("D19").Value = "3.35"
$ws.Range("E19").Value = "  +0.89%  "

$ws.Range("D20").Value = "0.0000228"
$ws.Range("E20").Value = "  +18.62%  "

$ws.Range("D21").Value = "13.28"
$ws.Range("E21").Value = "  -5.63%  "

$ws.Range("D22").Value = "426.96"
$ws.Range("E22").Value = "  -2.76%  "

$ws.Range("D23").Value = "8.40"
$ws.Range("E23").Value = "  -5.97%  "

$ws.Range("D24").Value = "4.91"
$ws.Range("E24").Value = "  -7.33%  "

$ws.Range("D25").Value = "5.42"
$ws.Range("E25").Value = "  +0.74%  "

$ws.Range("D26").Value = "11.53"
$ws.Range("E26").Value = "  -5.59%  "

$ws.Range("D27").Value = "80.05"
$ws.Range("E27").Value = "  +4.56%  "

$ws.Range("D28").Value = "3.280.72"
$ws.Range("E28").Value = "  -6.65%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  -12.85%  "

$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("D32").Value = "4.03"
$ws.Range("E32").Value = "  +9.56%  "

$ws.Range("D33").Value = "8.23"
$ws.Range("E33").Value = "  -6.72%  "

$ws.Range("D34").Value = "513.25"
$ws.Range("E34").Value = "  -8.78%  "

$ws.Range("D35").Value = "7.13"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").Value = "1.30"
$ws.Range("E36").Value = "  +0.63%  "

$ws.Range("D37").Value = "1.84"
$ws.Range("E37").Value = "  -4.69%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.136"
$ws.Range("E38").Value = "  +3.00%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "21.92"
$ws.Range("E39").Value = "  -3.59%  "

$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "22.23"
$ws.Range("E40").Value = "  -0.93%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("E42").Value = "  -0.17%  "

$ws.Range("D43").Value = "1.87"
$ws.Range("E43").Value = "  -6.92%  "

$ws.Range("D44").Value = "0.365"
$ws.Range("E44").Value = "  -8.02%  "

$ws.Range("D45").Value = "145.94"
$ws.Range("E45").Value = "  -3.06%  "

$ws.Range("D46").Value = "43.70"
$ws.Range("E46").Value = "  -0.08%  "

$ws.Range("E47").Value = "  -2.80%  "

$ws.Range("D48").Value = "166.09"
$ws.Range("E48").Value = "  -8.35%  "

$ws.Range("D49").Value = "0.727"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").Value = "24.58"
$ws.Range("E50").Value = "  -1.90%  "

$ws.Range("D51").Value = "1.19"
$ws.Range("E51").Value = "  -7.57%  "
